$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update "Pre Baseline Phase" (column B) values for each question row.
# These were placeholder values; the parent's initial measurements have now
# been filled in.
$ws.Range("B3").Value  = "Somewhat worse"
$ws.Range("B4").Value  = "Somewhat worse"
$ws.Range("B5").Value  = "Somewhat worse"
$ws.Range("B6").Value  = "Somewhat worse"
$ws.Range("B7").Value  = "A lot worse"
$ws.Range("B8").Value  = "Somewhat worse"
$ws.Range("B9").Value  = "Somewhat worse"
$ws.Range("B10").Value = "Somewhat worse"
$ws.Range("B11").Value = "A lot worse"
$ws.Range("B12").Value = "Somewhat worse"
$ws.Range("B13").Value = "Somewhat worse"
$ws.Range("B14").Value = "Somewhat worse"
$ws.Range("B15").Value = "A lot worse"
